# Updates the "cryptos" price/volume snapshot table on the active sheet.
# Numeric-looking price strings (column D) are prefixed with a leading
# apostrophe so Excel stores them as text (matching the source data's
# inlineStr cell type) instead of silently coercing them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.822.82'
$ws.Range('E2').Value = '  -1.39%  '
$ws.Range('D3').Value = '1.761.38'
$ws.Range('E3').Value = '  -2.98%  '
$ws.Range('D4').Value = "'1.002"
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = "'322.67"
$ws.Range('E5').Value = '  -2.10%  '
$ws.Range('D6').Value = "'1.001"
$ws.Range('E6').Value = '  +0.15%  '
$ws.Range('D7').Value = "'0.4264"
$ws.Range('E7').Value = '  -3.91%  '
$ws.Range('D8').Value = "'0.3624"
$ws.Range('E8').Value = '  -2.98%  '
$ws.Range('D9').Value = "'0.07580"
$ws.Range('E9').Value = '  -1.59%  '
$ws.Range('D10').Value = "'42.69"
$ws.Range('E10').Value = '  -4.66%  '
$ws.Range('D11').Value = "'1.095"
$ws.Range('E11').Value = '  -2.64%  '
$ws.Range('E12').Value = '  -0.03%  '
$ws.Range('D13').Value = "'20.73"
$ws.Range('E13').Value = '  -5.74%  '
$ws.Range('D14').Value = "'6.073"
$ws.Range('E14').Value = '  -3.54%  '
$ws.Range('D15').Value = "'7.277"
$ws.Range('E15').Value = '  -3.42%  '
$ws.Range('D16').Value = '1.759.39'
$ws.Range('E16').Value = '  -3.37%  '
$ws.Range('D17').Value = "'91.42"
$ws.Range('E17').Value = '  -2.67%  '
$ws.Range('D18').Value = "'0.00001067"
$ws.Range('E18').Value = '  -1.42%  '
$ws.Range('D19').Value = "'0.06381"
$ws.Range('E19').Value = '  -1.84%  '
$ws.Range('D20').Value = "'0.9998"
$ws.Range('E20').Value = '  +0.24%  '
$ws.Range('D21').Value = "'17.08"
$ws.Range('E21').Value = '  -2.46%  '
$ws.Range('D22').Value = "'5.914"
$ws.Range('E22').Value = '  -5.42%  '
$ws.Range('D23').Value = '27.864.16'
$ws.Range('E23').Value = '  -1.52%  '
$ws.Range('D24').Value = "'11.25"
$ws.Range('E24').Value = '  -3.94%  '
$ws.Range('D25').Value = "'2.121"
$ws.Range('E25').Value = '  +2.83%  '
$ws.Range('D26').Value = "'160.22"
$ws.Range('E26').Value = '  +2.90%  '
$ws.Range('D27').Value = "'20.29"
$ws.Range('E27').Value = '  -1.71%  '
$ws.Range('D28').Value = '1.967.96'
$ws.Range('E28').Value = '  -2.81%  '
$ws.Range('D29').Value = "'2.149"
$ws.Range('E29').Value = '  -7.58%  '
$ws.Range('D30').Value = "'124.91"
$ws.Range('E30').Value = '  -2.02%  '
$ws.Range('D31').Value = "'1.122"
$ws.Range('E31').Value = '  -6.50%  '
$ws.Range('D32').Value = "'3.680"
$ws.Range('E32').Value = '  +0.32%  '
$ws.Range('D33').Value = "'5.580"
$ws.Range('E33').Value = '  -4.65%  '
$ws.Range('D34').Value = "'0.08891"
$ws.Range('E34').Value = '  -3.68%  '
$ws.Range('D35').Value = "'12.26"
$ws.Range('E35').Value = '  -5.93%  '
$ws.Range('D36').Value = "'0.02302"
$ws.Range('E36').Value = '  -1.76%  '
$ws.Range('D37').Value = "'0.2110"
$ws.Range('E37').Value = '  -2.81%  '
$ws.Range('D38').Value = "'0.06036"
$ws.Range('E38').Value = '  -2.76%  '
$ws.Range('B39').Value = 'TheSandbox'
$ws.Range('C39').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D39').Value = "'0.6366"
$ws.Range('E39').Value = '  -3.04%  '
$ws.Range('B40').Value = 'InternetComputer(DFINITY)'
$ws.Range('C40').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D40').Value = "'4.987"
$ws.Range('E40').Value = '  -3.67%  '
$ws.Range('D41').Value = "'1.183"
$ws.Range('E41').Value = '  -0.99%  '
$ws.Range('D42').Value = "'0.9994"
$ws.Range('E42').Value = '  +0.19%  '
$ws.Range('D43').Value = "'7.917"
$ws.Range('E43').Value = '  -2.10%  '
$ws.Range('D44').Value = "'1.396"
$ws.Range('E44').Value = '  +0.57%  '
$ws.Range('D45').Value = "'13.30"
$ws.Range('E45').Value = '  -4.47%  '
$ws.Range('D46').Value = "'0.5884"
$ws.Range('E46').Value = '  -3.04%  '
$ws.Range('D47').Value = "'3.695"
$ws.Range('E47').Value = '  -1.70%  '
$ws.Range('E48').Value = '  -2.10%  '
$ws.Range('D49').Value = "'123.09"
$ws.Range('E49').Value = '  -2.79%  '
$ws.Range('D50').Value = "'1.184"
$ws.Range('E50').Value = '  +2.90%  '
$ws.Range('D51').Value = "'0.06835"
$ws.Range('E51').Value = '  -2.09%  '
